$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 878.6
$ws.Range("I19").Value = 933.6667
$ws.Range("K19").Value = 933.6667
$ws.Range("M19").Value = -758.6667
$ws.Range("H38").Value = 265.33334
$ws.Range("I38").Value = 48.18182
$ws.Range("K38").Value = 144.54546
$ws.Range("M38").Value = 227.45454
$ws.Range("H121").Value = 9150
$ws.Range("J121").Value = 9800
$ws.Range("L121").Value = 29400
$ws.Range("N121").Value = -32894
$ws.Range("H137").Value = 1670.7368
$ws.Range("I137").Value = 1392.7333
$ws.Range("J137").Value = 2713.25
$ws.Range("K137").Value = 4178.199900000001
$ws.Range("L137").Value = 8139.75
$ws.Range("M137").Value = -1628.199900000001
$ws.Range("N137").Value = -13239.75
$ws.Range("H141").Value = 3265.32
$ws.Range("I141").Value = 2675.4211
$ws.Range("K141").Value = 8026.263300000001
$ws.Range("M141").Value = -2846.263300000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3194.6667
$ws.Range("I2").Value = 2980.5715
$ws.Range("J2").Value = 3944
$ws.Range("K2").Value = 2980.5715
$ws.Range("L2").Value = 3944
$ws.Range("M2").Value = -2867.5715
$ws.Range("N2").Value = -4170
$ws.Range("H74").Value = 35715724
$ws.Range("I74").Value = 58824076
$ws.Range("J74").Value = 2813
$ws.Range("K74").Value = 58824076
$ws.Range("L74").Value = 2813
$ws.Range("M74").Value = -58823202
$ws.Range("N74").Value = -4561
$ws.Range("H77").Value = 35715724
$ws.Range("I77").Value = 58824076
$ws.Range("J77").Value = 2813
$ws.Range("K77").Value = 294120380
$ws.Range("L77").Value = 14065
$ws.Range("M77").Value = -294116012
$ws.Range("N77").Value = -22801
$ws.Range("H110").Value = 777.1667
$ws.Range("I110").Value = 670.75
$ws.Range("K110").Value = 670.75
$ws.Range("M110").Value = 1374.25
$ws.Range("H116").Value = 3194.6667
$ws.Range("I116").Value = 2980.5715
$ws.Range("J116").Value = 3944
$ws.Range("K116").Value = 2980.5715
$ws.Range("L116").Value = 3944
$ws.Range("M116").Value = -686.5715
$ws.Range("N116").Value = -8532

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3194.6667
$ws.Range("I3").Value = 2980.5715
$ws.Range("J3").Value = 3944
$ws.Range("K3").Value = 2980.5715
$ws.Range("L3").Value = 3944
$ws.Range("M3").Value = -2866.5715
$ws.Range("N3").Value = -4172
$ws.Range("H86").Value = 1512.878
$ws.Range("I86").Value = 1342.1482
$ws.Range("K86").Value = 1342.1482
$ws.Range("M86").Value = -219.1482000000001
$ws.Range("H89").Value = 1512.878
$ws.Range("I89").Value = 1342.1482
$ws.Range("K89").Value = 6710.741
$ws.Range("M89").Value = -1094.741
$ws.Range("H134").Value = 3156.9023
$ws.Range("I134").Value = 3557.6177
$ws.Range("J134").Value = 1210.5714
$ws.Range("K134").Value = 10672.8531
$ws.Range("L134").Value = 3631.7142
$ws.Range("M134").Value = -8137.8531
$ws.Range("N134").Value = -8701.7142

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16671.848
$ws.Range("I58").Value = 1301.381
$ws.Range("J58").Value = 43570.168
$ws.Range("K58").Value = 1301.381
$ws.Range("L58").Value = 43570.168
$ws.Range("M58").Value = -1098.381
$ws.Range("N58").Value = -43976.168
$ws.Range("H94").Value = 3439.7144
$ws.Range("I94").Value = 619
$ws.Range("K94").Value = 619
$ws.Range("M94").Value = -168
$ws.Range("H95").Value = 36000
$ws.Range("J95").Value = 36000
$ws.Range("L95").Value = 36000
$ws.Range("H99").Value = 20836616
$ws.Range("I99").Value = 3030.9092
$ws.Range("J99").Value = 38465036
$ws.Range("K99").Value = 3030.9092
$ws.Range("L99").Value = 38465036
$ws.Range("M99").Value = -1532.9092
$ws.Range("N99").Value = -38468032
$ws.Range("H126").Value = 20836616
$ws.Range("I126").Value = 3030.9092
$ws.Range("J126").Value = 38465036
$ws.Range("K126").Value = 9092.7276
$ws.Range("L126").Value = 115395108
$ws.Range("M126").Value = -6622.7276
$ws.Range("N126").Value = -115400048
$ws.Range("H136").Value = 16671.848
$ws.Range("I136").Value = 1301.381
$ws.Range("J136").Value = 43570.168
$ws.Range("K136").Value = 3904.143
$ws.Range("L136").Value = 130710.504
$ws.Range("M136").Value = -1354.143
$ws.Range("N136").Value = -135810.504

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 73.166664
$ws.Range("J12").Value = 95.333336
$ws.Range("L12").Value = 286.000008
$ws.Range("N12").Value = -632.000008
$ws.Range("H68").Value = 1238.75
$ws.Range("J68").Value = 1377.5
$ws.Range("L68").Value = 4132.5
$ws.Range("N68").Value = -5754.5
$ws.Range("H71").Value = 1238.75
$ws.Range("J71").Value = 1377.5
$ws.Range("L71").Value = 12397.5
$ws.Range("N71").Value = -20509.5
$ws.Range("H75").Value = 685
$ws.Range("I75").Value = 685
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 2055
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 685
$ws.Range("I78").Value = 685
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 6165
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H87").Value = 20059.176
$ws.Range("I87").Value = 13102.8
$ws.Range("K87").Value = 39308.39999999999
$ws.Range("M87").Value = -38060.39999999999
$ws.Range("H90").Value = 20059.176
$ws.Range("I90").Value = 13102.8
$ws.Range("K90").Value = 117925.2
$ws.Range("M90").Value = -111685.2
$ws.Range("H121").Value = 1025.3226
$ws.Range("I121").Value = 571.25
$ws.Range("J121").Value = 1092.5927
$ws.Range("K121").Value = 1713.75
$ws.Range("L121").Value = 3277.7781
$ws.Range("M121").Value = -403.75
$ws.Range("N121").Value = -5897.7781
$ws.Range("H131").Value = 797.46
$ws.Range("J131").Value = 795.375
$ws.Range("L131").Value = 2386.125
$ws.Range("N131").Value = -12466.125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 60000
$ws.Range("H70").Value = 8946400
$ws.Range("I70").Value = 5150
$ws.Range("J70").Value = 20868066
$ws.Range("K70").Value = 5150
$ws.Range("L70").Value = 20868066
$ws.Range("M70").Value = -4880
$ws.Range("N70").Value = -20868606
$ws.Range("H73").Value = 8946400
$ws.Range("I73").Value = 5150
$ws.Range("J73").Value = 20868066
$ws.Range("K73").Value = 5150
$ws.Range("L73").Value = 20868066
$ws.Range("M73").Value = -4214
$ws.Range("N73").Value = -20869938
$ws.Range("H80").Value = 3633.3333
$ws.Range("I80").Value = 2600
$ws.Range("J80").Value = 4253.3335
$ws.Range("K80").Value = 2600
$ws.Range("L80").Value = 4253.3335
$ws.Range("M80").Value = -1602
$ws.Range("N80").Value = -6249.3335
$ws.Range("H83").Value = 3633.3333
$ws.Range("I83").Value = 2600
$ws.Range("J83").Value = 4253.3335
$ws.Range("K83").Value = 13000
$ws.Range("L83").Value = 21266.6675
$ws.Range("M83").Value = -8008
$ws.Range("N83").Value = -31250.6675
$ws.Range("H97").Value = 1758.8148
$ws.Range("I97").Value = 1688.4117
$ws.Range("J97").Value = 1878.5
$ws.Range("K97").Value = 1688.4117
$ws.Range("L97").Value = 1878.5
$ws.Range("M97").Value = -1192.4117
$ws.Range("N97").Value = -2870.5
$ws.Range("H119").Value = 43753.332
$ws.Range("J119").Value = 43753.332
$ws.Range("L119").Value = 43753.332
$ws.Range("N119").Value = -53429.332

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3542.9092
$ws.Range("I40").Value = 2951.889
$ws.Range("J40").Value = 6202.5
$ws.Range("K40").Value = 2951.889
$ws.Range("L40").Value = 6202.5
$ws.Range("M40").Value = -2815.889
$ws.Range("N40").Value = -6474.5
$ws.Range("H82").Value = 3150.3333
$ws.Range("I82").Value = 3150.3333
$ws.Range("K82").Value = 3150.3333
$ws.Range("M82").Value = -2789.3333
$ws.Range("H85").Value = 3150.3333
$ws.Range("I85").Value = 3150.3333
$ws.Range("K85").Value = 3150.3333
$ws.Range("M85").Value = -1902.3333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 3000
$ws.Range("J21").Value = 3000
$ws.Range("L21").Value = 3000
$ws.Range("N21").Value = -3470
$ws.Range("H35").Value = 3000
$ws.Range("J35").Value = 3000
$ws.Range("L35").Value = 3000
$ws.Range("N35").Value = -3580
$ws.Range("H107").Value = 2674298.2
$ws.Range("I107").Value = 547.625
$ws.Range("J107").Value = 5050965.5
$ws.Range("K107").Value = 1642.875
$ws.Range("L107").Value = 15152896.5
$ws.Range("M107").Value = 277.125
$ws.Range("N107").Value = -15156736.5
$ws.Range("H122").Value = 845.2973
$ws.Range("I122").Value = 827.6818
$ws.Range("J122").Value = 871.13336
$ws.Range("K122").Value = 2483.0454
$ws.Range("L122").Value = 2613.40008
$ws.Range("M122").Value = -33.04539999999997
$ws.Range("N122").Value = -7513.40008
$ws.Range("H136").Value = 37040700
$ws.Range("I136").Value = 55557428
$ws.Range("K136").Value = 166672284
$ws.Range("M136").Value = -166669734
